$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 999
$ws.Range("I11").Value = 999
$ws.Range("K11").Value = 999
$ws.Range("M11").Value = -859

$ws.Range("H40").Value = 1553.4615
$ws.Range("I40").Value = 1122.9412
$ws.Range("J40").Value = 2366.6667
$ws.Range("K40").Value = 1122.9412
$ws.Range("L40").Value = 2366.6667
$ws.Range("M40").Value = -947.9412
$ws.Range("N40").Value = -2716.6667

$ws.Range("H51").Value = 2176
$ws.Range("I51").Value = 1931.25
$ws.Range("J51").Value = 2611.111
$ws.Range("K51").Value = 1931.25
$ws.Range("L51").Value = 2611.111
$ws.Range("M51").Value = -1447.25
$ws.Range("N51").Value = -3579.111

$ws.Range("H64").Value = 2927.6191
$ws.Range("I64").Value = 2913.3333
$ws.Range("J64").Value = 2963.3333
$ws.Range("K64").Value = 2913.3333
$ws.Range("L64").Value = 2963.3333
$ws.Range("M64").Value = -2665.3333
$ws.Range("N64").Value = -3459.3333

$ws.Range("H67").Value = 2927.6191
$ws.Range("I67").Value = 2913.3333
$ws.Range("J67").Value = 2963.3333
$ws.Range("K67").Value = 2913.3333
$ws.Range("L67").Value = 2963.3333
$ws.Range("M67").Value = -2055.3333
$ws.Range("N67").Value = -4679.3333

$ws.Range("H82").Value = 3861.111
$ws.Range("I82").Value = 1330
$ws.Range("J82").Value = 5471.8184
$ws.Range("K82").Value = 3990
$ws.Range("L82").Value = 16415.4552
$ws.Range("M82").Value = -3584
$ws.Range("N82").Value = -17227.4552

$ws.Range("H85").Value = 3861.111
$ws.Range("I85").Value = 1330
$ws.Range("J85").Value = 5471.8184
$ws.Range("K85").Value = 3990
$ws.Range("L85").Value = 16415.4552
$ws.Range("M85").Value = -2586
$ws.Range("N85").Value = -19223.4552

$ws.Range("H86").Value = 6201.304
$ws.Range("I86").Value = 17802.5
$ws.Range("J86").Value = 2106.7646
$ws.Range("K86").Value = 17802.5
$ws.Range("L86").Value = 2106.7646
$ws.Range("M86").Value = -16679.5
$ws.Range("N86").Value = -4352.7646

$ws.Range("H89").Value = 6201.304
$ws.Range("I89").Value = 17802.5
$ws.Range("J89").Value = 2106.7646
$ws.Range("K89").Value = 89012.5
$ws.Range("L89").Value = 10533.823
$ws.Range("M89").Value = -83396.5
$ws.Range("N89").Value = -21765.823

$ws.Range("H100").Value = 3154.9656
$ws.Range("I100").Value = 3231
$ws.Range("J100").Value = 3061.3845
$ws.Range("K100").Value = 3231
$ws.Range("L100").Value = 3061.3845
$ws.Range("M100").Value = -2690
$ws.Range("N100").Value = -4143.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6759208.5
$ws.Range("I2").Value = 19232334
$ws.Range("K2").Value = 19232334
$ws.Range("M2").Value = -19232221

$ws.Range("H74").Value = 959.8570999999999
$ws.Range("I74").Value = 892.4737
$ws.Range("J74").Value = 1600
$ws.Range("K74").Value = 892.4737
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = -18.47370000000001
$ws.Range("N74").Value = -3348

$ws.Range("H77").Value = 959.8570999999999
$ws.Range("I77").Value = 892.4737
$ws.Range("J77").Value = 1600
$ws.Range("K77").Value = 4462.3685
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = -94.36850000000049
$ws.Range("N77").Value = -16736

$ws.Range("H116").Value = 6759208.5
$ws.Range("I116").Value = 19232334
$ws.Range("K116").Value = 19232334
$ws.Range("M116").Value = -19230040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6759208.5
$ws.Range("I3").Value = 19232334
$ws.Range("K3").Value = 19232334
$ws.Range("M3").Value = -19232220

$ws.Range("H86").Value = 2003.2
$ws.Range("I86").Value = 1298.6666
$ws.Range("J86").Value = 3060
$ws.Range("K86").Value = 1298.6666
$ws.Range("L86").Value = 3060
$ws.Range("M86").Value = -175.6666
$ws.Range("N86").Value = -5306

$ws.Range("H89").Value = 2003.2
$ws.Range("I89").Value = 1298.6666
$ws.Range("J89").Value = 3060
$ws.Range("K89").Value = 6493.333000000001
$ws.Range("L89").Value = 15300
$ws.Range("M89").Value = -877.3330000000005
$ws.Range("N89").Value = -26532

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3338438.5
$ws.Range("I31").Value = 5884371
$ws.Range("K31").Value = 5884371
$ws.Range("M31").Value = -5884076

$ws.Range("H34").Value = 3338438.5
$ws.Range("I34").Value = 5884371
$ws.Range("K34").Value = 5884371
$ws.Range("M34").Value = -5884169

$ws.Range("H132").Value = 3285.516
$ws.Range("I132").Value = 1935.8667
$ws.Range("J132").Value = 4550.8125
$ws.Range("K132").Value = 5807.6001
$ws.Range("L132").Value = 13652.4375
$ws.Range("M132").Value = -3277.6001
$ws.Range("N132").Value = -18712.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 559.8
$ws.Range("I13").Value = 449.75
$ws.Range("K13").Value = 1349.25
$ws.Range("M13").Value = -1181.25

$ws.Range("H139").Value = 7581324
$ws.Range("I139").Value = 11907054
$ws.Range("J139").Value = 11296.333
$ws.Range("K139").Value = 35721162
$ws.Range("L139").Value = 33888.999
$ws.Range("M139").Value = -35716022
$ws.Range("N139").Value = -44168.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1016179.4
$ws.Range("I3").Value = 1250223.2
$ws.Range("K3").Value = 1250223.2
$ws.Range("M3").Value = -1250107.2

$ws.Range("H9").Value = 16281.333
$ws.Range("I9").Value = 930.8570999999999
$ws.Range("K9").Value = 930.8570999999999
$ws.Range("M9").Value = -760.8570999999999

$ws.Range("H92").Value = 15166.333
$ws.Range("J92").Value = 15166.333
$ws.Range("L92").Value = 15166.333
$ws.Range("N92").Value = -18910.333

$ws.Range("H97").Value = 1298.55
$ws.Range("I97").Value = 1041.8182
$ws.Range("J97").Value = 1612.3334
$ws.Range("K97").Value = 1041.8182
$ws.Range("L97").Value = 1612.3334
$ws.Range("M97").Value = -545.8181999999999
$ws.Range("N97").Value = -2604.3334

$ws.Range("H122").Value = 3736.0476
$ws.Range("I122").Value = 4600.875
$ws.Range("J122").Value = 3203.8462
$ws.Range("K122").Value = 13802.625
$ws.Range("L122").Value = 9611.5386
$ws.Range("M122").Value = -11352.625
$ws.Range("N122").Value = -14511.5386

$ws.Range("H132").Value = 3607.4614
$ws.Range("I132").Value = 1985.5714
$ws.Range("K132").Value = 5956.7142
$ws.Range("M132").Value = -3426.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4600
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 5900
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 5900
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -6124

$ws.Range("H12").Value = 10003
$ws.Range("I12").Value = 10003
$ws.Range("K12").Value = 10003
$ws.Range("M12").Value = -9833

$ws.Range("H126").Value = 4600
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5900
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 17700
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -22640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 581.0833
$ws.Range("I100").Value = 543
$ws.Range("K100").Value = 1086
$ws.Range("M100").Value = -545

$ws.Range("H122").Value = 668327.9399999999
$ws.Range("I122").Value = 1667712.4
$ws.Range("K122").Value = 5003137.199999999
$ws.Range("M122").Value = -5000687.199999999

$ws.Range("H126").Value = 9092846
$ws.Range("I126").Value = 959.4
$ws.Range("K126").Value = 2878.2
$ws.Range("M126").Value = -408.1999999999998
